$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.008
$ws.Range("C5").Value = 0.02
$ws.Range("C10").Value = 0.595
$ws.Range("C11").Value = 0.965
$ws.Range("C12").Value = 0.13
$ws.Range("C13").Value = 0.474
